$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.028014397248625755
$ws.Range("C2").Value = 0.013654418289661407
$ws.Range("D2").Value = 0.008439145050942898
$ws.Range("E2").Value = 0.006082390900701284
$ws.Range("F2").Value = 0.0002644160413183272
$ws.Range("J2").Value = 0.12751227617263794
$ws.Range("K2").Value = 1.4451524019241333
